# Apply the authored edit:
#  1. Refresh the cached "datetimeFigureOut" footer field (Insert > Header &
#     Footer > Date and time > Update automatically) from 17-Sep-25 to
#     29-Oct-25 on the slide master and every slide layout.
#  2. Remove the last two slides (slide 8 and slide 9) from the deck.

$p = $ppt.ActivePresentation

$oldDate = "17-Sep-25"
$newDate = "29-Oct-25"

$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        $isDatePlaceholder = $false
        if ($sh.Type -eq 14) {
            # msoPlaceholder
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $isDatePlaceholder = $true
            }
        } elseif ($sh.Name -like "Date Placeholder*") {
            $isDatePlaceholder = $true
        }
        if ($isDatePlaceholder -and $sh.TextFrame.HasText -and $sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master footer date field.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's footer date field.
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    Update-DatePlaceholder $layout.Shapes
}

# Drop the trailing two slides (formerly slide8.xml / slide9.xml).
while ($p.Slides.Count -gt 7) {
    $p.Slides.Item($p.Slides.Count).Delete()
}
